$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "canonical SMILES" column header (D2)
$ws.Range("D2").Value = "canonical SMILES"

# Populate the new column with the canonical SMILES values (duplicating the
# canonical isomeric SMILES already present in column C for each molecule)
$ws.Range("D3").Value = $ws.Range("C3").Value2
$ws.Range("D4").Value = $ws.Range("C4").Value2
$ws.Range("D5").Value = $ws.Range("C5").Value2
$ws.Range("D6").Value = $ws.Range("C6").Value2

# Match the column width used in the target workbook
$ws.Columns.Item(4).ColumnWidth = 36
